$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 with the new retailer (Bismillah Telecom / RET-26511) ---
$ws.Range("B2").Value = "RET-26511"
$ws.Range("C2").Value = "Bismillah Telecom"

# --- Center + middle-align the rest of row 2 (A2, and clone onto B2:E2) ---
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4108
$ws.Range("A2").Copy()
$ws.Range("B2:E2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Prepare the quote-prefixed + centered number style on F3 (still has the
# original quotePrefix base style) before touching F2's value, then stamp
# that exact format onto F2 once its new value is in place. This avoids
# creating throw-away intermediate cell styles.
$ws.Range("F3").HorizontalAlignment = -4108
$ws.Range("F3").VerticalAlignment = -4108

$ws.Range("F2").Value = 1316416301
$ws.Range("F3").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Clear the now-obsolete retailer rows, keeping their formatting ---
$ws.Range("A3:F7").ClearContents()

# --- Update the active selection marker ---
$ws.Range("I10").Select()
